$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "rsync this script to your day8 directory." paragraph:
#    add "You can " prefix, and extend the sentence with the github mention.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "rsync this script to your day8 directory.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "You can rsync this script to your day8 directory or pull the script from github under day8 scripts.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "...ntasks and local cores equal to 34." -> "...ntasks AND local cores..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "ntasks and local cores equal to 34.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "ntasks AND local cores equal to 34.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Reposition the three floating pictures (Top/Left, in points).
#    EMU targets: 5415915 / 6076569 / (-843280, 7753858)
#    NOTE: Shapes.Item(N) setters resolve N against raw document order, so
#    the indices below were determined empirically (document order), not
#    the z-order used when enumerating/printing the collection.
# ---------------------------------------------------------------------------
$d.Shapes.Item(5).Top = 5415915 / 12700.0
$d.Shapes.Item(6).Top = 6076569 / 12700.0
$d.Shapes.Item(7).Top = 7753858 / 12700.0
$d.Shapes.Item(7).Left = -843280 / 12700.0

# ---------------------------------------------------------------------------
# 4) "Path to the directory with the fastq's:" -> "Path to the fastq's directory:"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Path to the directory with the fastq’s:",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Path to the fastq’s directory:",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 5) "...were downloaded using the following curl command which can be found
#    on the 10x website:" -> "...were downloaded ahead of time using the
#    following curl command, also from the 10x website:"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "were downloaded using the following curl command which can be found on the 10x website:",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "were downloaded ahead of time using the following curl command, also from the 10x website:",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 6) Delete the (now obsolete) "*MAY*DELETE*" paragraph block: the blank
#    paragraph right after the curl sentence, the *MAY*DELETE* paragraph,
#    the blank paragraph after it, and the "export PATH=..." paragraph.
#    The following blank paragraph (before "Now, run the sbatch script")
#    is kept.
# ---------------------------------------------------------------------------
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*MAY*DELETE*") {
        $startPara = $i - 1
        $endPara = $i + 2
    }
}
if ($startPara -ne $null) {
    $rStart = $d.Paragraphs($startPara).Range.Start
    $rEnd = $d.Paragraphs($endPara).Range.End
    $d.Range($rStart, $rEnd).Delete() | Out-Null
}

Write-Output "done"
